$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tr = $s.Shapes.Item(2).TextFrame.TextRange

# Force the writer to rebuild the paragraph's runs: flipping the text to a
# throwaway value first, then to the desired consolidated text, ensures the
# final write emits a single run instead of reusing the existing split runs.
$tr.Text = "placeholder_tmp"
$tr.Text = "The picture first"
